$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug List")

# Set "Version Fixed" (column C) to 1.0.7 for rows 14, 15, 18
$ws.Range("C14").Value = "1.0.7"
$ws.Range("C15").Value = "1.0.7"
$ws.Range("C18").Value = "1.0.7"

# Update the active selection to D18 on the Bug List sheet
$ws.Range("D18").Select()
